# Word COM-interop script: convert the inline DrawingML picture in temp2.docx
# into a legacy VML <w:pict> shape, and relocate the "_GoBack" bookmark from
# around the picture run to sit right after the "Template2" text run.

$d = $word.ActiveDocument

$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'
$vNs = 'urn:schemas-microsoft-com:vml'
$oNs = 'urn:schemas-microsoft-com:office:office'
$rNs = 'http://schemas.openxmlformats.org/officeDocument/2006/relationships'

# --- 1. Add the "_GoBack" bookmark right after the "Template2" run -----------
# (paragraph 1 = "Template2"). Replacing the paragraph's whole text range with
# itself plus the bookmark tags keeps the bookmark from spanning the
# paragraph mark (a bookmark collapsed exactly on the mark gets pushed into
# the following paragraph instead).
$p1 = $d.Paragraphs(1).Range
$p1Xml = "<w:p xmlns:w='$wNs' w:rsidR='001E11A2' w:rsidRDefault='00F42264'>" +
         "<w:pPr><w:rPr><w:lang w:val='en-GB'/></w:rPr></w:pPr>" +
         "<w:r><w:rPr><w:lang w:val='en-GB'/></w:rPr><w:t>Template2</w:t></w:r>" +
         "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
         "<w:bookmarkEnd w:id='0'/>" +
         "</w:p>"
$d.Range($p1.Start, $p1.End - 1).InsertXML($p1Xml) | Out-Null

# --- 2. Swap the DrawingML <w:drawing> picture for a VML <w:pict> shape -----
# Targeting the whole paragraph 3 range (pPr + run + the old bookmark pair
# around it) replaces the old bookmarkStart/bookmarkEnd as well, matching the
# diff (they are simply dropped, not carried over).
$p3 = $d.Paragraphs(3).Range
$pictXml =
  "<w:p xmlns:w='$wNs' xmlns:v='$vNs' xmlns:o='$oNs' xmlns:r='$rNs' " +
  "w:rsidR='003C1847' w:rsidRPr='00EF6969' w:rsidRDefault='003C1847'>" +
    "<w:pPr><w:rPr><w:lang w:val='en-GB'/></w:rPr></w:pPr>" +
    "<w:r>" +
      "<w:rPr><w:noProof/><w:lang w:eastAsia='de-DE'/></w:rPr>" +
      "<w:pict>" +
        "<v:shapetype id='_x0000_t75' coordsize='21600,21600' o:spt='75' " +
          "o:preferrelative='t' path='m@4@5l@4@11@9@11@9@5xe' filled='f' stroked='f'>" +
          "<v:stroke joinstyle='miter'/>" +
          "<v:formulas>" +
            "<v:f eqn='if lineDrawn pixelLineWidth 0'/>" +
            "<v:f eqn='sum @0 1 0'/>" +
            "<v:f eqn='sum 0 0 @1'/>" +
            "<v:f eqn='prod @2 1 2'/>" +
            "<v:f eqn='prod @3 21600 pixelWidth'/>" +
            "<v:f eqn='prod @3 21600 pixelHeight'/>" +
            "<v:f eqn='sum @0 0 1'/>" +
            "<v:f eqn='prod @6 1 2'/>" +
            "<v:f eqn='prod @7 21600 pixelWidth'/>" +
            "<v:f eqn='sum @8 21600 0'/>" +
            "<v:f eqn='prod @7 21600 pixelHeight'/>" +
            "<v:f eqn='sum @10 21600 0'/>" +
          "</v:formulas>" +
          "<v:path o:extrusionok='f' gradientshapeok='t' o:connecttype='rect'/>" +
          "<o:lock v:ext='edit' aspectratio='t'/>" +
        "</v:shapetype>" +
        "<v:shape id='_x0000_i1025' type='#_x0000_t75' style='width:453pt;height:255pt'>" +
          "<v:imagedata r:id='rId4' o:title='relax'/>" +
        "</v:shape>" +
      "</w:pict>" +
    "</w:r>" +
  "</w:p>"
$d.Range($p3.Start, $p3.End).InsertXML($pictXml) | Out-Null

Write-Output "edit applied"
